$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LSTM Model Summary")
$ws.Range("B73").Value = "epochs"
Write-Host "done"
